$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder columns: "Throttle" moves from column B to column D,
# shifting "Elevator" and "Rudder" one column to the left.
$ws.Range("B1").Value = "Elevator"
$ws.Range("C1").Value = "Rudder"
$ws.Range("D1").Value = "Throttle"

# Header font style change: Calibri -> Helvetica
$ws.Range("A1:D1").Font.Name = "Helvetica"

# Row height tweak for the header row
$ws.Rows("1").RowHeight = 18

# Update the active selection to D1 (single cell)
$ws.Range("D1").Select()
